$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Estadisticos 1P')
$ws.Range('A2').Value = 'Ingles II'
$ws.Range('B2').Value = '2AEV'
$ws.Range('C2').Value = 36
$ws.Range('D2').Value = 0
$ws.Range('E2').Value = 17
$ws.Range('F2').Value = 19
$ws.Range('G2').Value = 52.78
$ws.Range('H2').Value = 5.7
$ws.Range('A3').Value = 'Ingles IV'
$ws.Range('B3').Value = '4AEM'
$ws.Range('C3').Value = 24
$ws.Range('D3').Value = 0
$ws.Range('E3').Value = 1
$ws.Range('F3').Value = 23
$ws.Range('G3').Value = 95.83
$ws.Range('H3').Value = 8.4
$ws.Range('A4').Value = 'Ingles IV'
$ws.Range('B4').Value = '4ALCM'
$ws.Range('C4').Value = 34
$ws.Range('D4').Value = 0
$ws.Range('E4').Value = 2
$ws.Range('F4').Value = 32
$ws.Range('G4').Value = 94.12
$ws.Range('H4').Value = 8.4
$ws.Range('A5').Value = 'Ingles IV'
$ws.Range('B5').Value = '4APM'
$ws.Range('C5').Value = 31
$ws.Range('D5').Value = 0
$ws.Range('E5').Value = 2
$ws.Range('F5').Value = 29
$ws.Range('G5').Value = 93.55
$ws.Range('H5').Value = 7.5
$ws.Range('A6').Value = 'Formación socioemocional IV'
$ws.Range('B6').Value = '4ARHM'
$ws.Range('C6').Value = 0
$ws.Range('D6').Value = 0
$ws.Range('E6').Value = 0
$ws.Range('F6').Value = 0
$ws.Range('G6').ClearContents()
$ws.Range('H6').ClearContents()
$ws.Range('A7').Value = 'Ingles IV'
$ws.Range('B7').Value = '4BEM'
$ws.Range('C7').Value = 24
$ws.Range('D7').Value = 0
$ws.Range('E7').Value = 1
$ws.Range('F7').Value = 23
$ws.Range('G7').Value = 95.83
$ws.Range('H7').Value = 8.6

$ws = $wb.Worksheets.Item('Estadisticos 2P')
$ws.Range('A2').Value = 'Ingles II'
$ws.Range('B2').Value = '2AEV'
$ws.Range('C2').Value = 36
$ws.Range('D2').Value = 0
$ws.Range('E2').Value = 10
$ws.Range('F2').Value = 26
$ws.Range('G2').Value = 72.22
$ws.Range('H2').Value = 5.7
$ws.Range('A3').Value = 'Ingles IV'
$ws.Range('B3').Value = '4AEM'
$ws.Range('C3').Value = 24
$ws.Range('D3').Value = 0
$ws.Range('E3').Value = 3
$ws.Range('F3').Value = 21
$ws.Range('G3').Value = 87.5
$ws.Range('H3').Value = 8.4
$ws.Range('A4').Value = 'Ingles IV'
$ws.Range('B4').Value = '4ALCM'
$ws.Range('C4').Value = 34
$ws.Range('D4').Value = 0
$ws.Range('E4').Value = 2
$ws.Range('F4').Value = 32
$ws.Range('G4').Value = 94.12
$ws.Range('H4').Value = 8.4
$ws.Range('A5').Value = 'Ingles IV'
$ws.Range('B5').Value = '4APM'
$ws.Range('C5').Value = 31
$ws.Range('D5').Value = 0
$ws.Range('E5').Value = 2
$ws.Range('F5').Value = 29
$ws.Range('G5').Value = 93.55
$ws.Range('H5').Value = 7.5
$ws.Range('A6').Value = 'Formación socioemocional IV'
$ws.Range('B6').Value = '4ARHM'
$ws.Range('C6').Value = 0
$ws.Range('D6').Value = 0
$ws.Range('E6').Value = 0
$ws.Range('F6').Value = 0
$ws.Range('G6').ClearContents()
$ws.Range('H6').ClearContents()
$ws.Range('A7').Value = 'Ingles IV'
$ws.Range('B7').Value = '4BEM'
$ws.Range('C7').Value = 24
$ws.Range('D7').Value = 0
$ws.Range('E7').Value = 0
$ws.Range('F7').Value = 24
$ws.Range('G7').Value = 100
$ws.Range('H7').Value = 8.6

$ws = $wb.Worksheets.Item('Estadisticos Final')
$ws.Range('A2').Value = 'Ingles II'
$ws.Range('B2').Value = '2AEV'
$ws.Range('C2').Value = 36
$ws.Range('D2').Value = 0
$ws.Range('E2').Value = 10
$ws.Range('F2').Value = 26
$ws.Range('G2').Value = 72.22
$ws.Range('H2').Value = 6.9
$ws.Range('A3').Value = 'Ingles IV'
$ws.Range('B3').Value = '4AEM'
$ws.Range('C3').Value = 24
$ws.Range('D3').Value = 0
$ws.Range('E3').Value = 3
$ws.Range('F3').Value = 21
$ws.Range('G3').Value = 87.5
$ws.Range('H3').Value = 8
$ws.Range('A4').Value = 'Ingles IV'
$ws.Range('B4').Value = '4ALCM'
$ws.Range('C4').Value = 34
$ws.Range('D4').Value = 0
$ws.Range('E4').Value = 2
$ws.Range('F4').Value = 32
$ws.Range('G4').Value = 94.12
$ws.Range('H4').Value = 8.9
$ws.Range('A5').Value = 'Ingles IV'
$ws.Range('B5').Value = '4APM'
$ws.Range('C5').Value = 31
$ws.Range('D5').Value = 0
$ws.Range('E5').Value = 2
$ws.Range('F5').Value = 29
$ws.Range('G5').Value = 93.55
$ws.Range('H5').Value = 7.9
$ws.Range('A6').Value = 'Formación socioemocional IV'
$ws.Range('B6').Value = '4ARHM'
$ws.Range('C6').Value = 0
$ws.Range('D6').Value = 0
$ws.Range('E6').Value = 0
$ws.Range('F6').Value = 0
$ws.Range('G6').ClearContents()
$ws.Range('H6').ClearContents()
$ws.Range('A7').Value = 'Ingles IV'
$ws.Range('B7').Value = '4BEM'
$ws.Range('C7').Value = 24
$ws.Range('D7').Value = 0
$ws.Range('E7').Value = 0
$ws.Range('F7').Value = 24
$ws.Range('G7').Value = 100
$ws.Range('H7').Value = 8.5

$ws = $wb.Worksheets.Item('Rescatables')
$ws.Range('A2').Value = 24330051920304
$ws.Range('B2').Value = 'ARMAS'
$ws.Range('C2').Value = 'SALINAS'
$ws.Range('D2').Value = 'JOSE GUSTAVO'
$ws.Range('E2').Value = 'Ingles II'
$ws.Range('F2').Value = '2AEV'
$ws.Range('G2').Value = 4
$ws.Range('A3').Value = 24330051920305
$ws.Range('B3').Value = 'MORALES'
$ws.Range('C3').Value = 'CUAHUA'
$ws.Range('D3').Value = 'ANDRES'
$ws.Range('E3').Value = 'Ingles II'
$ws.Range('F3').Value = '2AEV'
$ws.Range('G3').Value = 4
$ws.Range('A4').Value = 24330051920113
$ws.Range('B4').Value = 'RAMOS'
$ws.Range('C4').Value = 'DE LA CRUZ'
$ws.Range('D4').Value = 'DEREK'
$ws.Range('E4').Value = 'Ingles II'
$ws.Range('F4').Value = '2AEV'
$ws.Range('G4').Value = 4
$ws.Range('A5').Value = 23330051920313
$ws.Range('B5').Value = 'VIVANCO'
$ws.Range('C5').Value = 'VIVANCO'
$ws.Range('D5').Value = 'LUIS AARON'
$ws.Range('E5').Value = 'Ingles IV'
$ws.Range('F5').Value = '4APM'
$ws.Range('G5').Value = 4
$ws.Range('A6').Value = 24330051920093
$ws.Range('B6').Value = 'ARIAS'
$ws.Range('C6').Value = 'SARMIENTO'
$ws.Range('D6').Value = 'URIEL ARTURO'
$ws.Range('E6').Value = 'Ingles II'
$ws.Range('F6').Value = '2AEV'
$ws.Range('G6').Value = 3
$ws.Range('A7').Value = 24330051920098
$ws.Range('B7').Value = 'CHICO'
$ws.Range('C7').Value = 'BALDERAS'
$ws.Range('D7').Value = 'YARETH'
$ws.Range('E7').Value = 'Ingles II'
$ws.Range('F7').Value = '2AEV'
$ws.Range('G7').Value = 3
$ws.Range('A8').Value = 24330051920144
$ws.Range('B8').Value = 'MUÑOZ'
$ws.Range('C8').Value = 'CORONA'
$ws.Range('D8').Value = 'JOSE ABEL'
$ws.Range('E8').Value = 'Ingles II'
$ws.Range('F8').Value = '2AEV'
$ws.Range('G8').Value = 3
$ws.Range('A9').Value = 24330051920143
$ws.Range('B9').Value = 'ROSAS'
$ws.Range('C9').Value = 'MEZA'
$ws.Range('D9').Value = 'CARLOS ANTONIO'
$ws.Range('E9').Value = 'Ingles II'
$ws.Range('F9').Value = '2AEV'
$ws.Range('G9').Value = 3
$ws.Range('A10').Value = 23330051920081
$ws.Range('B10').Value = 'CARRERA'
$ws.Range('C10').Value = 'MOLINA'
$ws.Range('D10').Value = 'MARIA DEL CARMEN'
$ws.Range('E10').Value = 'Ingles IV'
$ws.Range('F10').Value = '4ALCM'
$ws.Range('G10').Value = 3
$ws.Range('A11').Value = 24330051920392
$ws.Range('B11').Value = 'CERON'
$ws.Range('C11').Value = 'GONZALEZ'
$ws.Range('D11').Value = 'LEVI SANTIAGO'
$ws.Range('E11').Value = 'Ingles II'
$ws.Range('F11').Value = '2AEV'
$ws.Range('G11').Value = 2
$ws.Range('A12').Value = 23330051920224
$ws.Range('B12').Value = 'DORANTES'
$ws.Range('C12').Value = 'PORRAS'
$ws.Range('D12').Value = 'ROBERTO'
$ws.Range('E12').Value = 'Ingles II'
$ws.Range('F12').Value = '2AEV'
$ws.Range('G12').Value = 2
$ws.Range('A13').Value = 23330051920113
$ws.Range('B13').Value = 'VASQUEZ'
$ws.Range('C13').Value = 'BAEZ'
$ws.Range('D13').Value = 'YAMILETH'
$ws.Range('E13').Value = 'Ingles IV'
$ws.Range('F13').Value = '4ALCM'
$ws.Range('G13').Value = 2
$ws.Range('A14').Value = 23330051920155
$ws.Range('B14').Value = 'CRUZ'
$ws.Range('C14').Value = 'NIEVES'
$ws.Range('D14').Value = 'ESTRELLA ESMERALDA'
$ws.Range('E14').Value = 'Ingles IV'
$ws.Range('F14').Value = '4APM'
$ws.Range('G14').Value = 1
